# Revert "Add snRNAseq-10xGenomics-v2 to the scrnaseq assays"
#
# The "assay_type list" sheet gained a row ("snRNAseq-10xGenomics-v2") in the
# original commit; this reverts it by removing that row again. Deleting the
# row shifts the subsequent entries up, shrinking the list from 7 to 6
# values, so the data-validation list on the "Export as TSV" sheet (column L,
# "assay_type") that points at that range also needs to be narrowed to match.

$wb = $excel.ActiveWorkbook

$listSheet = $wb.Worksheets.Item("assay_type list")

# Find & remove the "snRNAseq-10xGenomics-v2" row from the list sheet.
$used = $listSheet.UsedRange
$rowCount = $used.Rows.Count
for ($i = 1; $i -le $rowCount; $i++) {
    $cell = $listSheet.Cells.Item($i, 1)
    if ($cell.Value2 -eq "snRNAseq-10xGenomics-v2") {
        $listSheet.Rows($i).Delete()
        break
    }
}

# The assay_type list shrank from 7 rows to 6 rows; keep the data validation
# on the main sheet's assay_type column (L) referencing the correct range.
$mainSheet = $wb.Worksheets.Item("Export as TSV")
$validation = $mainSheet.Range("L2:L1048576").Validation
$validation.Formula1 = "='assay_type list'!`$A`$1:`$A`$6"
